$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D edits to stay as text (avoid Excel auto-numeric coercion),
# then clear the temporary format so styling matches the original (no explicit style).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.463.08'
$ws.Range("E2").Value = '  +0.04%  '

$ws.Range("D3").Value = '1.825.16'
$ws.Range("E3").Value = '  -0.22%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '314.57'
$ws.Range("E5").Value = '  -0.83%  '

$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").Value = '0.5093'
$ws.Range("E7").Value = '  -4.75%  '

$ws.Range("D8").Value = '0.3928'
$ws.Range("E8").Value = '  -2.57%  '

$ws.Range("D9").Value = '0.07680'
$ws.Range("E9").Value = '  +1.16%  '

$ws.Range("D10").Value = '41.89'
$ws.Range("E10").Value = '  +0.05%  '

$ws.Range("D11").Value = '1.110'
$ws.Range("E11").Value = '  +0.21%  '

$ws.Range("D12").Value = '21.02'
$ws.Range("E12").Value = '  +0.58%  '

$ws.Range("D13").Value = '6.271'
$ws.Range("E13").Value = '  -1.05%  '

$ws.Range("D14").Value = '1.001'
$ws.Range("E14").Value = '  -0.02%  '

$ws.Range("D15").Value = '7.527'
$ws.Range("E15").Value = '  -0.32%  '

$ws.Range("D16").Value = '1.826.43'
$ws.Range("E16").Value = '  +0.77%  '

$ws.Range("D17").Value = '93.01'
$ws.Range("E17").Value = '  +4.03%  '

$ws.Range("D18").Value = '0.00001105'
$ws.Range("E18").Value = '  +3.20%  '

$ws.Range("D19").Value = '0.06642'
$ws.Range("E19").Value = '  +0.46%  '

$ws.Range("D20").Value = '17.78'
$ws.Range("E20").Value = '  +1.07%  '

$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("D22").Value = '6.118'
$ws.Range("E22").Value = '  +0.85%  '

$ws.Range("D23").Value = '28.473.70'
$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("D24").Value = '11.23'
$ws.Range("E24").Value = '  -0.57%  '

$ws.Range("D25").Value = '2.253'
$ws.Range("E25").Value = '  +5.04%  '

$ws.Range("D26").Value = '21.54'
$ws.Range("E26").Value = '  +4.69%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '155.98'
$ws.Range("E27").Value = '  -0.25%  '

$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.034.17'
$ws.Range("E28").Value = '  -0.26%  '

$ws.Range("D29").Value = '2.403'
$ws.Range("E29").Value = '  -3.66%  '

$ws.Range("D30").Value = '124.61'
$ws.Range("E30").Value = '  +0.91%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '1.111'
$ws.Range("E31").Value = '  -0.71%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '0.1099'
$ws.Range("E32").Value = '  +0.87%  '

$ws.Range("D33").Value = '5.665'
$ws.Range("E33").Value = '  -0.16%  '

$ws.Range("D34").Value = '3.659'
$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").Value = '0.07048'
$ws.Range("E35").Value = '  -1.65%  '

$ws.Range("D36").Value = '0.2214'
$ws.Range("E36").Value = '  -2.11%  '

$ws.Range("D37").Value = '0.02326'
$ws.Range("E37").Value = '  -0.72%  '

$ws.Range("D38").Value = '5.173'
$ws.Range("E38").Value = '  -0.98%  '

$ws.Range("D39").Value = '8.769'
$ws.Range("E39").Value = '  -0.85%  '

$ws.Range("D40").Value = '0.6270'
$ws.Range("E40").Value = '  -0.16%  '

$ws.Range("D41").Value = '11.20'
$ws.Range("E41").Value = '  -1.07%  '

$ws.Range("D42").Value = '1.175'
$ws.Range("E42").Value = '  -0.45%  '

$ws.Range("E43").Value = '  -0.02%  '

$ws.Range("E44").Value = '  -0.77%  '

$ws.Range("D45").Value = '13.45'
$ws.Range("E45").Value = '  +0.28%  '

$ws.Range("D46").Value = '3.731'
$ws.Range("E46").Value = '  +0.73%  '

$ws.Range("D47").Value = '0.5884'
$ws.Range("E47").Value = '  +0.42%  '

$ws.Range("D48").Value = '124.27'
$ws.Range("E48").Value = '  -1.35%  '

$ws.Range("D49").Value = '1.983'
$ws.Range("E49").Value = '  -0.54%  '

$ws.Range("D50").Value = '1.193'
$ws.Range("E50").Value = '  -0.11%  '

$ws.Range("E51").Value = '  +0.08%  '

# Restore default (General) formatting on column D so styles match the source file.
$ws.Range("D2:D51").ClearFormats()
